# Mangrove Dec-24 Profit and Loss — monthly roll-forward update.
# Renames the two "...June'24" distributor commission/campaign line items
# to "...Dec'24" (this workbook is reused every month) and fills in the
# amounts received for this month's Distributor GA & SAF Commission and
# the Accelerate income line, then scrolls/selects near the bottom of the
# sheet (where the report now ends) as the last thing the user did before
# saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Distributor GA & SAF Commission — relabel June'24 -> Dec'24 and
# record this month's amount.
$ws.Range("A13").Value = "Distributor GA & SAF Commission Dec'24"
$ws.Range("C13").Value = 2475

# Row 14: Distributor Campaign Shera partner — relabel June'24 -> Dec'24.
$ws.Range("A14").Value = "Distributor Campaign Shera partner Dec'24"

# Row 15: Accelerate — label unchanged, but this month's amount is filled in.
$ws.Range("C15").Value = 2736

# C17 (Total income) and F43 (Profit/Loss) are formulas and recalculate
# automatically from the above.

# Leave the sheet scrolled/selected near the bottom, matching where the
# user finished working.
$ws.Range("F52").Select()
